$d = $word.ActiveDocument

$d.Content.Find.Execute("2023-05-02 Tuesday", $true, $true, $false, $false, $false, $true, 1, $false, "2023-05-03 Wednesday", 2) | Out-Null
$d.Content.Find.Execute("98-17=", $true, $true, $false, $false, $false, $true, 1, $false, "96+1=", 2) | Out-Null
$d.Content.Find.Execute("3+46=", $true, $true, $false, $false, $false, $true, 1, $false, "35+60=", 2) | Out-Null
$d.Content.Find.Execute("37-11=", $true, $true, $false, $false, $false, $true, 1, $false, "46-24=", 2) | Out-Null
$d.Content.Find.Execute("79-66=", $true, $true, $false, $false, $false, $true, 1, $false, "19+9=", 2) | Out-Null
$d.Content.Find.Execute("21+29=", $true, $true, $false, $false, $false, $true, 1, $false, "41-25=", 2) | Out-Null
$d.Content.Find.Execute("78-8=", $true, $true, $false, $false, $false, $true, 1, $false, "55+12=", 2) | Out-Null
$d.Content.Find.Execute("54-10=", $true, $true, $false, $false, $false, $true, 1, $false, "12+66=", 2) | Out-Null
$d.Content.Find.Execute("99-73=", $true, $true, $false, $false, $false, $true, 1, $false, "58-45=", 2) | Out-Null
$d.Content.Find.Execute("46+31=", $true, $true, $false, $false, $false, $true, 1, $false, "20+37=", 2) | Out-Null
$d.Content.Find.Execute("0+34=", $true, $true, $false, $false, $false, $true, 1, $false, "11+59=", 2) | Out-Null
$d.Content.Find.Execute("38+43=", $true, $true, $false, $false, $false, $true, 1, $false, "93-43=", 2) | Out-Null
$d.Content.Find.Execute("41-17=", $true, $true, $false, $false, $false, $true, 1, $false, "72-71=", 2) | Out-Null
$d.Content.Find.Execute("60-44=", $true, $true, $false, $false, $false, $true, 1, $false, "70+21=", 2) | Out-Null
$d.Content.Find.Execute("46+18=", $true, $true, $false, $false, $false, $true, 1, $false, "56+27=", 2) | Out-Null
$d.Content.Find.Execute("68-35=", $true, $true, $false, $false, $false, $true, 1, $false, "38-1=", 2) | Out-Null
$d.Content.Find.Execute("34-27=", $true, $true, $false, $false, $false, $true, 1, $false, "88+6=", 2) | Out-Null
$d.Content.Find.Execute("45-23=", $true, $true, $false, $false, $false, $true, 1, $false, "32+22=", 2) | Out-Null
$d.Content.Find.Execute("59+26=", $true, $true, $false, $false, $false, $true, 1, $false, "67-61=", 2) | Out-Null
$d.Content.Find.Execute("59-29=", $true, $true, $false, $false, $false, $true, 1, $false, "56+12=", 2) | Out-Null
$d.Content.Find.Execute("18+18=", $true, $true, $false, $false, $false, $true, 1, $false, "5+63=", 2) | Out-Null
$d.Content.Find.Execute("26+64=", $true, $true, $false, $false, $false, $true, 1, $false, "14+53=", 2) | Out-Null
$d.Content.Find.Execute("62-47=", $true, $true, $false, $false, $false, $true, 1, $false, "11+39=", 2) | Out-Null
$d.Content.Find.Execute("71-43=", $true, $true, $false, $false, $false, $true, 1, $false, "54-20=", 2) | Out-Null
$d.Content.Find.Execute("21+64=", $true, $true, $false, $false, $false, $true, 1, $false, "41+19=", 2) | Out-Null
$d.Content.Find.Execute("4+8=", $true, $true, $false, $false, $false, $true, 1, $false, "61-12=", 2) | Out-Null
$d.Content.Find.Execute("27+54=", $true, $true, $false, $false, $false, $true, 1, $false, "87-47=", 2) | Out-Null
$d.Content.Find.Execute("35-27=", $true, $true, $false, $false, $false, $true, 1, $false, "6+31=", 2) | Out-Null
$d.Content.Find.Execute("22+0=", $true, $true, $false, $false, $false, $true, 1, $false, "87-46=", 2) | Out-Null
$d.Content.Find.Execute("76+9=", $true, $true, $false, $false, $false, $true, 1, $false, "67-24=", 2) | Out-Null
$d.Content.Find.Execute("19+14=", $true, $true, $false, $false, $false, $true, 1, $false, "28+6=", 2) | Out-Null
$d.Content.Find.Execute("63-12=", $true, $true, $false, $false, $false, $true, 1, $false, "27-7=", 2) | Out-Null
$d.Content.Find.Execute("55-15=", $true, $true, $false, $false, $false, $true, 1, $false, "34-24=", 2) | Out-Null
$d.Content.Find.Execute("62-5=", $true, $true, $false, $false, $false, $true, 1, $false, "62-9=", 2) | Out-Null
$d.Content.Find.Execute("90-15=", $true, $true, $false, $false, $false, $true, 1, $false, "2+34=", 2) | Out-Null
$d.Content.Find.Execute("44-37=", $true, $true, $false, $false, $false, $true, 1, $false, "49-20=", 2) | Out-Null
$d.Content.Find.Execute("28+61=", $true, $true, $false, $false, $false, $true, 1, $false, "57-17=", 2) | Out-Null
$d.Content.Find.Execute("44+18=", $true, $true, $false, $false, $false, $true, 1, $false, "2+89=", 2) | Out-Null
$d.Content.Find.Execute("71-52=", $true, $true, $false, $false, $false, $true, 1, $false, "22+46=", 2) | Out-Null
$d.Content.Find.Execute("19+2=", $true, $true, $false, $false, $false, $true, 1, $false, "15+52=", 2) | Out-Null
$d.Content.Find.Execute("53-42=", $true, $true, $false, $false, $false, $true, 1, $false, "80-67=", 2) | Out-Null
$d.Content.Find.Execute("58+40=", $true, $true, $false, $false, $false, $true, 1, $false, "77-30=", 2) | Out-Null
$d.Content.Find.Execute("58-37=", $true, $true, $false, $false, $false, $true, 1, $false, "80-27=", 2) | Out-Null
$d.Content.Find.Execute("2+94=", $true, $true, $false, $false, $false, $true, 1, $false, "16-0=", 2) | Out-Null
$d.Content.Find.Execute("99-7=", $true, $true, $false, $false, $false, $true, 1, $false, "14+79=", 2) | Out-Null
$d.Content.Find.Execute("28+40=", $true, $true, $false, $false, $false, $true, 1, $false, "64-39=", 2) | Out-Null
$d.Content.Find.Execute("60+36=", $true, $true, $false, $false, $false, $true, 1, $false, "17+24=", 2) | Out-Null
$d.Content.Find.Execute("6+25=", $true, $true, $false, $false, $false, $true, 1, $false, "13+9=", 2) | Out-Null
$d.Content.Find.Execute("71-16=", $true, $true, $false, $false, $false, $true, 1, $false, "54+34=", 2) | Out-Null
$d.Content.Find.Execute("78-46=", $true, $true, $false, $false, $false, $true, 1, $false, "51-29=", 2) | Out-Null
$d.Content.Find.Execute("50+0=", $true, $true, $false, $false, $false, $true, 1, $false, "36-29=", 2) | Out-Null
$d.Content.Find.Execute("54-17=", $true, $true, $false, $false, $false, $true, 1, $false, "48-43=", 2) | Out-Null
$d.Content.Find.Execute("84-22=", $true, $true, $false, $false, $false, $true, 1, $false, "67-54=", 2) | Out-Null
$d.Content.Find.Execute("38+25=", $true, $true, $false, $false, $false, $true, 1, $false, "42+32=", 2) | Out-Null
$d.Content.Find.Execute("11+69=", $true, $true, $false, $false, $false, $true, 1, $false, "13+58=", 2) | Out-Null
$d.Content.Find.Execute("47-21=", $true, $true, $false, $false, $false, $true, 1, $false, "80-67=", 2) | Out-Null
$d.Content.Find.Execute("52+28=", $true, $true, $false, $false, $false, $true, 1, $false, "99-34=", 2) | Out-Null
$d.Content.Find.Execute("24+63=", $true, $true, $false, $false, $false, $true, 1, $false, "98-82=", 2) | Out-Null
$d.Content.Find.Execute("68-46=", $true, $true, $false, $false, $false, $true, 1, $false, "91-11=", 2) | Out-Null
$d.Content.Find.Execute("58+16=", $true, $true, $false, $false, $false, $true, 1, $false, "27+21=", 2) | Out-Null
$d.Content.Find.Execute("27+36=", $true, $true, $false, $false, $false, $true, 1, $false, "54-52=", 2) | Out-Null
$d.Content.Find.Execute("99-83=", $true, $true, $false, $false, $false, $true, 1, $false, "10+53=", 2) | Out-Null
$d.Content.Find.Execute("13+18=", $true, $true, $false, $false, $false, $true, 1, $false, "33-16=", 2) | Out-Null
$d.Content.Find.Execute("38+19=", $true, $true, $false, $false, $false, $true, 1, $false, "36+0=", 2) | Out-Null
$d.Content.Find.Execute("16+46=", $true, $true, $false, $false, $false, $true, 1, $false, "17+26=", 2) | Out-Null
$d.Content.Find.Execute("77-61=", $true, $true, $false, $false, $false, $true, 1, $false, "24-3=", 2) | Out-Null
$d.Content.Find.Execute("73-26=", $true, $true, $false, $false, $false, $true, 1, $false, "72-39=", 2) | Out-Null
$d.Content.Find.Execute("5+9=", $true, $true, $false, $false, $false, $true, 1, $false, "46-24=", 2) | Out-Null
$d.Content.Find.Execute("39-12=", $true, $true, $false, $false, $false, $true, 1, $false, "63+29=", 2) | Out-Null
$d.Content.Find.Execute("21+11=", $true, $true, $false, $false, $false, $true, 1, $false, "53-40=", 2) | Out-Null
$d.Content.Find.Execute("63+15=", $true, $true, $false, $false, $false, $true, 1, $false, "35+23=", 2) | Out-Null
$d.Content.Find.Execute("49+1=", $true, $true, $false, $false, $false, $true, 1, $false, "20+20=", 2) | Out-Null
$d.Content.Find.Execute("73-63=", $true, $true, $false, $false, $false, $true, 1, $false, "30+51=", 2) | Out-Null
$d.Content.Find.Execute("30+26=", $true, $true, $false, $false, $false, $true, 1, $false, "31+7=", 2) | Out-Null
$d.Content.Find.Execute("76-52=", $true, $true, $false, $false, $false, $true, 1, $false, "32+55=", 2) | Out-Null
$d.Content.Find.Execute("38+58=", $true, $true, $false, $false, $false, $true, 1, $false, "19+67=", 2) | Out-Null
$d.Content.Find.Execute("62+10=", $true, $true, $false, $false, $false, $true, 1, $false, "33+22=", 2) | Out-Null
$d.Content.Find.Execute("56+17=", $true, $true, $false, $false, $false, $true, 1, $false, "32+65=", 2) | Out-Null
$d.Content.Find.Execute("35+39=", $true, $true, $false, $false, $false, $true, 1, $false, "87-3=", 2) | Out-Null
$d.Content.Find.Execute("4-0=", $true, $true, $false, $false, $false, $true, 1, $false, "97-0=", 2) | Out-Null
$d.Content.Find.Execute("5-4=", $true, $true, $false, $false, $false, $true, 1, $false, "81-13=", 2) | Out-Null
$d.Content.Find.Execute("99-50=", $true, $true, $false, $false, $false, $true, 1, $false, "2+50=", 2) | Out-Null
$d.Content.Find.Execute("23-22=", $true, $true, $false, $false, $false, $true, 1, $false, "46-16=", 2) | Out-Null
$d.Content.Find.Execute("12+73=", $true, $true, $false, $false, $false, $true, 1, $false, "20+40=", 2) | Out-Null
$d.Content.Find.Execute("82-79=", $true, $true, $false, $false, $false, $true, 1, $false, "7+65=", 2) | Out-Null
$d.Content.Find.Execute("24-11=", $true, $true, $false, $false, $false, $true, 1, $false, "44-7=", 2) | Out-Null
$d.Content.Find.Execute("34-32=", $true, $true, $false, $false, $false, $true, 1, $false, "77-36=", 2) | Out-Null
$d.Content.Find.Execute("36+4=", $true, $true, $false, $false, $false, $true, 1, $false, "45+16=", 2) | Out-Null
$d.Content.Find.Execute("90-26=", $true, $true, $false, $false, $false, $true, 1, $false, "57+32=", 2) | Out-Null
$d.Content.Find.Execute("96-19=", $true, $true, $false, $false, $false, $true, 1, $false, "50-25=", 2) | Out-Null
$d.Content.Find.Execute("84-60=", $true, $true, $false, $false, $false, $true, 1, $false, "82-34=", 2) | Out-Null
$d.Content.Find.Execute("72-45=", $true, $true, $false, $false, $false, $true, 1, $false, "68-33=", 2) | Out-Null
$d.Content.Find.Execute("16+45=", $true, $true, $false, $false, $false, $true, 1, $false, "97-66=", 2) | Out-Null
$d.Content.Find.Execute("52-22=", $true, $true, $false, $false, $false, $true, 1, $false, "67-45=", 2) | Out-Null
$d.Content.Find.Execute("31+65=", $true, $true, $false, $false, $false, $true, 1, $false, "79+7=", 2) | Out-Null
$d.Content.Find.Execute("79-28=", $true, $true, $false, $false, $false, $true, 1, $false, "87+1=", 2) | Out-Null
$d.Content.Find.Execute("20+68=", $true, $true, $false, $false, $false, $true, 1, $false, "40+30=", 2) | Out-Null
$d.Content.Find.Execute("3+50=", $true, $true, $false, $false, $false, $true, 1, $false, "36-2=", 2) | Out-Null
$d.Content.Find.Execute("80-32=", $true, $true, $false, $false, $false, $true, 1, $false, "66-60=", 2) | Out-Null
$d.Content.Find.Execute("44-23=", $true, $true, $false, $false, $false, $true, 1, $false, "88-49=", 2) | Out-Null
$d.Content.Find.Execute("63+17=", $true, $true, $false, $false, $false, $true, 1, $false, "2+61=", 2) | Out-Null
